# Implementation of web service.
# - Bump the cached "datetimeFigureOut" date placeholder text from
#   19-10-2022 to 20-10-2022 on the slide master and every slide layout
#   (the file was re-saved a day later, so PowerPoint re-cached the
#   auto date field everywhere it appears).
# - Point the hello-service endpoint URL at its new location.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder ("Date Placeholder N") cached text: 19-10-2022 -> 20-10-2022
#    Shape index of the date placeholder varies per layout, so address
#    each one explicitly (verified against each layout's shape order).
# ---------------------------------------------------------------------
$newDate = "20-10-2022"

$master = $p.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

$layouts = $master.CustomLayouts
$layouts.Item(1).Shapes.Item(3).TextFrame.TextRange.Text = $newDate    # Title Slide
$layouts.Item(2).Shapes.Item(3).TextFrame.TextRange.Text = $newDate    # Title and Content
$layouts.Item(3).Shapes.Item(3).TextFrame.TextRange.Text = $newDate    # Section Header
$layouts.Item(4).Shapes.Item(4).TextFrame.TextRange.Text = $newDate    # Two Content
$layouts.Item(5).Shapes.Item(6).TextFrame.TextRange.Text = $newDate    # Comparison
$layouts.Item(6).Shapes.Item(2).TextFrame.TextRange.Text = $newDate    # Title Only
$layouts.Item(7).Shapes.Item(1).TextFrame.TextRange.Text = $newDate    # Blank
$layouts.Item(8).Shapes.Item(4).TextFrame.TextRange.Text = $newDate    # Content with Caption
$layouts.Item(9).Shapes.Item(4).TextFrame.TextRange.Text = $newDate    # Picture with Caption
$layouts.Item(10).Shapes.Item(3).TextFrame.TextRange.Text = $newDate   # Title and Vertical Text
$layouts.Item(11).Shapes.Item(3).TextFrame.TextRange.Text = $newDate   # Vertical Title and Text

# ---------------------------------------------------------------------
# 2. Endpoint URL text on slide 5 ("Problem 1 Phase 1 - Design"):
#    http://helloservice/v1/user/id -> http://app/hello/v1/user/id
#    Re-target just that paragraph. The link text shares a prefix/suffix
#    with the old value, so first stamp an unrelated placeholder value to
#    force a clean single run, then set the real text - this keeps the
#    hyperlink run properties intact instead of fragmenting into several
#    runs that share the same formatting.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(5)
$endpointShape = $slide.Shapes.Item(2)
$urlPara = $endpointShape.TextFrame.TextRange.Paragraphs(2, 1)
$urlPara.Text = "_"
$urlPara = $endpointShape.TextFrame.TextRange.Paragraphs(2, 1)
$urlPara.Text = "http://app/hello/v1/user/id"
